# Fruta / hortaliza, semanal
# Feria Lagunitas de Puerto Montt - Frutilla: add a new weekly price update
# (2022-04-05, serial 44656) for "Primera" and "Segunda" quality, inserted
# right before the existing 2021-11-15 / 2022-01-31 rows (which shift down
# from rows 225-226 to rows 227-228).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert two new rows starting at row 225. This pushes the
# existing rows 225-226 down to rows 227-228, leaving two blank rows at
# 225-226 for the new data.
$ws.Range("A225:A226").EntireRow.Insert()

# --- Row 225: brand-new row (date 2022-04-05, "Primera" quality) ---
$ws.Cells.Item(225, 1).Value = 4
$ws.Cells.Item(225, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(225, 3).Value = "Los Lagos"
$ws.Cells.Item(225, 4).Value = 44656        # D225 Fecha
$ws.Cells.Item(225, 5).Value = 10
$ws.Cells.Item(225, 6).Value = "Fruta"
$ws.Cells.Item(225, 7).Value = 100101
$ws.Cells.Item(225, 8).Value = "Berries"
$ws.Cells.Item(225, 9).Value = 100112025
$ws.Cells.Item(225, 10).Value = "Frutilla"
$ws.Cells.Item(225, 11).Value = "Sin especificar"
$ws.Cells.Item(225, 12).Value = "Primera"
$ws.Cells.Item(225, 13).Value = 400         # M225 Volumen
$ws.Cells.Item(225, 14).Value = 8000        # N225 Precio minimo
$ws.Cells.Item(225, 15).Value = 8500        # O225 Precio maximo
$ws.Cells.Item(225, 16).Value = 8250        # P225 Precio promedio ponderado
$ws.Cells.Item(225, 17).Value = "$/caja 7 kilos"          # Q225 Unidad de comercializacion
$ws.Cells.Item(225, 18).Value = "Región de La Araucanía"  # R225 Origen
$ws.Cells.Item(225, 19).Value = 1179        # S225 Precio $/Kg
$ws.Cells.Item(225, 20).Value = 7           # T225 Kg / unidad

# --- Row 226: brand-new row (same date, "Segunda" quality) ---
$ws.Cells.Item(226, 1).Value = 4
$ws.Cells.Item(226, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(226, 3).Value = "Los Lagos"
$ws.Cells.Item(226, 4).Value = 44656
$ws.Cells.Item(226, 5).Value = 10
$ws.Cells.Item(226, 6).Value = "Fruta"
$ws.Cells.Item(226, 7).Value = 100101
$ws.Cells.Item(226, 8).Value = "Berries"
$ws.Cells.Item(226, 9).Value = 100112025
$ws.Cells.Item(226, 10).Value = "Frutilla"
$ws.Cells.Item(226, 11).Value = "Sin especificar"
$ws.Cells.Item(226, 12).Value = "Segunda"
$ws.Cells.Item(226, 13).Value = 200
$ws.Cells.Item(226, 14).Value = 6000
$ws.Cells.Item(226, 15).Value = 6000
$ws.Cells.Item(226, 16).Value = 6000
$ws.Cells.Item(226, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(226, 18).Value = "Región de La Araucanía"
$ws.Cells.Item(226, 19).Value = 857
$ws.Cells.Item(226, 20).Value = 7

# Rows 227 and 228 already hold the former rows 225/226 content (shifted
# down automatically by the insert above), so nothing else to do there.
